$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 109
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2719
$ws.Range("H33").Value = 1180
$ws.Range("I33").Value = 645.125
$ws.Range("J33").Value = 2249.75
$ws.Range("K33").Value = 645.125
$ws.Range("L33").Value = 2249.75
$ws.Range("M33").Value = -416.125
$ws.Range("N33").Value = -2707.75
$ws.Range("H38").Value = 1820.6666
$ws.Range("J38").Value = 2100
$ws.Range("L38").Value = 6300
$ws.Range("N38").Value = -7044
$ws.Range("H55").Value = 327.41666
$ws.Range("I55").Value = 112.833336
$ws.Range("J55").Value = 542
$ws.Range("K55").Value = 112.833336
$ws.Range("L55").Value = 542
$ws.Range("M55").Value = 101.166664
$ws.Range("N55").Value = -970
$ws.Range("H58").Value = 2042.1818
$ws.Range("I58").Value = 669.8889
$ws.Range("J58").Value = 2992.2307
$ws.Range("K58").Value = 2009.6667
$ws.Range("L58").Value = 8976.6921
$ws.Range("M58").Value = -1859.6667
$ws.Range("N58").Value = -9276.6921
$ws.Range("H80").Value = 63844.97
$ws.Range("I80").Value = 88518.39
$ws.Range("J80").Value = 790.6667
$ws.Range("K80").Value = 265555.17
$ws.Range("L80").Value = 2372.0001
$ws.Range("M80").Value = -264557.17
$ws.Range("N80").Value = -4368.0001
$ws.Range("H83").Value = 63844.97
$ws.Range("I83").Value = 88518.39
$ws.Range("J83").Value = 790.6667
$ws.Range("K83").Value = 796665.51
$ws.Range("L83").Value = 7116.0003
$ws.Range("M83").Value = -791673.51
$ws.Range("N83").Value = -17100.0003
$ws.Range("H107").Value = 6519.913
$ws.Range("I107").Value = 8874.267
$ws.Range("K107").Value = 8874.267
$ws.Range("M107").Value = -6954.267
$ws.Range("H111").Value = 4000
$ws.Range("I111").Value = 4000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 12000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -8933
$ws.Range("H112").Value = 3780.4443
$ws.Range("J112").Value = 5000
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17216
$ws.Range("H132").Value = 7158.0415
$ws.Range("I132").Value = 7491.421
$ws.Range("K132").Value = 22474.263
$ws.Range("M132").Value = -19944.263
$ws.Range("H137").Value = 18193
$ws.Range("I137").Value = 22450.9
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 67352.70000000001
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -64802.70000000001
$ws.Range("N137").Value = -17100
$ws.Range("N9").ClearContents()
$ws.Range("N111").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 803549.9
$ws.Range("I122").Value = 3082.2058
$ws.Range("J122").Value = 2099545.2
$ws.Range("K122").Value = 9246.617400000001
$ws.Range("L122").Value = 6298635.600000001
$ws.Range("M122").Value = -6796.617400000001
$ws.Range("N122").Value = -6303535.600000001
$ws.Range("H132").Value = 4272.162
$ws.Range("I132").Value = 3356.7083
$ws.Range("J132").Value = 5962.231
$ws.Range("K132").Value = 10070.1249
$ws.Range("L132").Value = 17886.693
$ws.Range("M132").Value = -7540.124899999999
$ws.Range("N132").Value = -22946.693

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1063.9131
$ws.Range("I22").Value = 684
$ws.Range("J22").Value = 1412.1666
$ws.Range("K22").Value = 684
$ws.Range("L22").Value = 1412.1666
$ws.Range("M22").Value = -334
$ws.Range("N22").Value = -2112.1666
$ws.Range("H99").Value = 5535108
$ws.Range("I99").Value = 10560229
$ws.Range("K99").Value = 10560229
$ws.Range("M99").Value = -10558731
$ws.Range("H107").Value = 8420.344999999999
$ws.Range("I107").Value = 11963.25
$ws.Range("J107").Value = 547.2222
$ws.Range("K107").Value = 11963.25
$ws.Range("L107").Value = 547.2222
$ws.Range("M107").Value = -10043.25
$ws.Range("N107").Value = -4387.2222
$ws.Range("H122").Value = 8647.9375
$ws.Range("I122").Value = 15474
$ws.Range("K122").Value = 46422
$ws.Range("M122").Value = -43972
$ws.Range("H126").Value = 5535108
$ws.Range("I126").Value = 10560229
$ws.Range("K126").Value = 31680687
$ws.Range("M126").Value = -31678217
$ws.Range("H133").Value = 80001
$ws.Range("J133").Value = 80001
$ws.Range("L133").Value = 80001
$ws.Range("N133").Value = -85061

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -2888
$ws.Range("N3").Value = -3224
$ws.Range("H23").Value = 5747524.5
$ws.Range("J23").Value = 8333721.5
$ws.Range("L23").Value = 25001164.5
$ws.Range("N23").Value = -25001634.5
$ws.Range("H34").Value = 1340.4
$ws.Range("I34").Value = 1338
$ws.Range("J34").Value = 1350
$ws.Range("K34").Value = 4014
$ws.Range("L34").Value = 4050
$ws.Range("M34").Value = -3930
$ws.Range("N34").Value = -4218
$ws.Range("H48").Value = 4062.8
$ws.Range("J48").Value = 4605
$ws.Range("L48").Value = 13815
$ws.Range("N48").Value = -14315
$ws.Range("H126").Value = 17226.8
$ws.Range("I126").Value = 2734.5
$ws.Range("J126").Value = 26888.334
$ws.Range("K126").Value = 8203.5
$ws.Range("L126").Value = 80665.00199999999
$ws.Range("M126").Value = -3263.5
$ws.Range("N126").Value = -90545.00199999999
$ws.Range("H133").Value = 12499
$ws.Range("J133").Value = 19999
$ws.Range("L133").Value = 59997
$ws.Range("N133").Value = -70117
$ws.Range("H134").Value = 2618.25
$ws.Range("I134").Value = 2618.25
$ws.Range("K134").Value = 7854.75
$ws.Range("M134").Value = -2784.75
$ws.Range("H136").Value = 2312
$ws.Range("I136").Value = 1543.3
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 4629.9
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = 470.1000000000004
$ws.Range("N136").Value = -40197
$ws.Range("H138").Value = 1275
$ws.Range("I138").Value = 830
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 2490
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = 2650
$ws.Range("N138").Value = -20780
$ws.Range("H139").Value = 3078636
$ws.Range("I139").Value = 4001027
$ws.Range("K139").Value = 12003081
$ws.Range("M139").Value = -11997941

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 10011667
$ws.Range("I14").Value = 10011667
$ws.Range("K14").Value = 10011667
$ws.Range("M14").Value = -10011499
$ws.Range("H97").Value = 5321.25
$ws.Range("I97").Value = 6486.0527
$ws.Range("K97").Value = 6486.0527
$ws.Range("M97").Value = -5990.0527
$ws.Range("H113").Value = 25340
$ws.Range("J113").Value = 4175
$ws.Range("L113").Value = 4175
$ws.Range("N113").Value = -8515
$ws.Range("H122").Value = 6591.25
$ws.Range("I122").Value = 4095.4243
$ws.Range("K122").Value = 12286.2729
$ws.Range("M122").Value = -9836.2729
$ws.Range("H132").Value = 2107.111
$ws.Range("I132").Value = 2120.5
$ws.Range("K132").Value = 6361.5
$ws.Range("M132").Value = -3831.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12994.333
$ws.Range("I20").Value = 12994.5
$ws.Range("K20").Value = 12994.5
$ws.Range("M20").Value = -12768.5
$ws.Range("H40").Value = 28089.523
$ws.Range("I40").Value = 42840.7
$ws.Range("K40").Value = 42840.7
$ws.Range("M40").Value = -42704.7
$ws.Range("H46").Value = 3533.0667
$ws.Range("I46").Value = 699.75
$ws.Range("K46").Value = 699.75
$ws.Range("M46").Value = -511.75
$ws.Range("H61").Value = 5483.25
$ws.Range("I61").Value = 3556.6843
$ws.Range("J61").Value = 12804.2
$ws.Range("K61").Value = 3556.6843
$ws.Range("L61").Value = 12804.2
$ws.Range("M61").Value = -3354.6843
$ws.Range("N61").Value = -13208.2
$ws.Range("H100").Value = 6124.875
$ws.Range("I100").Value = 3250.5
$ws.Range("J100").Value = 8999.25
$ws.Range("K100").Value = 3250.5
$ws.Range("L100").Value = 8999.25
$ws.Range("M100").Value = -2709.5
$ws.Range("N100").Value = -10081.25
$ws.Range("H113").Value = 5483.25
$ws.Range("I113").Value = 3556.6843
$ws.Range("J113").Value = 12804.2
$ws.Range("K113").Value = 3556.6843
$ws.Range("L113").Value = 12804.2
$ws.Range("M113").Value = -1386.6843
$ws.Range("N113").Value = -17144.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9653
$ws.Range("I81").Value = 11259.9
$ws.Range("J81").Value = 4296.6665
$ws.Range("K81").Value = 22519.8
$ws.Range("L81").Value = 8593.333000000001
$ws.Range("M81").Value = -21458.8
$ws.Range("N81").Value = -10715.333
$ws.Range("H84").Value = 9653
$ws.Range("I84").Value = 11259.9
$ws.Range("J84").Value = 4296.6665
$ws.Range("K84").Value = 112599
$ws.Range("L84").Value = 42966.665
$ws.Range("M84").Value = -107295
$ws.Range("N84").Value = -53574.665
$ws.Range("H107").Value = 9780.102999999999
$ws.Range("J107").Value = 50549.832
$ws.Range("L107").Value = 151649.496
$ws.Range("N107").Value = -155489.496
$ws.Range("H122").Value = 5445.12
$ws.Range("I122").Value = 2793.7273
$ws.Range("K122").Value = 8381.1819
$ws.Range("M122").Value = -5931.1819
$ws.Range("H138").Value = 94000
$ws.Range("J138").Value = 94000
$ws.Range("L138").Value = 94000
$ws.Range("N138").Value = -104280
